$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.96670711040496826
$ws.Range("D1").Value = 4.6447315216064453
$ws.Range("F1").Value = 4.6447315216064453

$ws.Range("A3").Value = 0.93554514646530151
$ws.Range("B3").Value = 0.77921265363693237
$ws.Range("C3").Value = 5.0259566307067871
$ws.Range("D3").Value = 3.6013355255126949
$ws.Range("E3").Value = 2.3132750988006592
$ws.Range("F3").Value = 2.028793573379517

$ws.Range("A7").Value = 1.3971821069717409
$ws.Range("B7").Value = 0.89828085899353027
$ws.Range("C7").Value = 2.6527543067932129
$ws.Range("D7").Value = 5.126030445098877
$ws.Range("E7").Value = 1.339913129806519
$ws.Range("F7").Value = 2.2339012622833252

$ws.Range("A17").Value = 1.0224477052688601
$ws.Range("C17").Value = 5.5734972953796387
$ws.Range("E17").Value = 2.3733184337615971
